$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group 1: style "old17" (green) -> J5:J10
$ws.Range("J7").Copy()
$dest1 = $excel.Union($ws.Range("J5"), $ws.Range("J6"), $ws.Range("J7"), $ws.Range("J8"), $ws.Range("J9"), $ws.Range("J10"))
$dest1.PasteSpecial(-4122)

# Group 2: style "old18" (gray) -> J11,J12,J14,J15
$ws.Range("J11").Copy()
$dest2 = $excel.Union($ws.Range("J11"), $ws.Range("J12"), $ws.Range("J14"), $ws.Range("J15"))
$dest2.PasteSpecial(-4122)

# Group 3: style "old25" (yellow) -> J13,J16
$ws.Range("J13").Copy()
$dest3 = $excel.Union($ws.Range("J13"), $ws.Range("J16"))
$dest3.PasteSpecial(-4122)
